$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row for 70b6c8ce-... (row 2) moves from "Ready for handoff"
# to "Handed back: in sync with en-US" now that it has been handed back.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Helper data per-language sheet: source file md + its handoff xlf info, plus
# the urls used for the new "Latest Target File" / "Latest Handback File"
# hyperlinks (same targets as the existing source / handoff-file links).
# ---------------------------------------------------------------------------

# zh-cn sheet -----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"

$wsZh.Range("F2").Value = "70b6c8ce-f45c-48c3-b209-a11b72d043a5.md"
$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/b236d27f6dd833f7f4a5e383568f97018d1d5d8a/e2e/70b6c8ce-f45c-48c3-b209-a11b72d043a5.md", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.md")

$wsZh.Range("G2").Value = "70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.zh-cn.xlf"
$wsZh.Range("G2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7571533ed125ae01df10791eb33567d14f922765/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.zh-cn.xlf", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.zh-cn.xlf")

$wsZh.Range("H2").Value = "2016-03-18 00:29:22"

# de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"

$wsDe.Range("F2").Value = "70b6c8ce-f45c-48c3-b209-a11b72d043a5.md"
$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/b236d27f6dd833f7f4a5e383568f97018d1d5d8a/e2e/70b6c8ce-f45c-48c3-b209-a11b72d043a5.md", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.md")

$wsDe.Range("G2").Value = "70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.de-de.xlf"
$wsDe.Range("G2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bad623942a2c732acb0ef0cc15d4a1bb697b0384/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.de-de.xlf", "", "", "70b6c8ce-f45c-48c3-b209-a11b72d043a5.263f2c9b563870b33ccaf31ef11bc6af2bef5be1.de-de.xlf")

$wsDe.Range("H2").Value = "2016-03-18 00:29:27"

Write-Host "Handback report generated."
